$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 16:22"
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7142426
$ws.Range("C4").Value = 2873
$ws.Range("D4").Value = 4400342
$ws.Range("E4").Value = 2535461
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 206623

$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 5752178
$ws.Range("C5").Value = 21994
$ws.Range("D5").Value = 4694860
$ws.Range("E5").Value = 965950
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 195
$ws.Range("H5").Value = 91368

$ws.Range("A32").Value = "Catar"
$ws.Range("B32").Value = 124425
$ws.Range("C32").Value = 250
$ws.Range("D32").Value = 121263
$ws.Range("E32").Value = 2950
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 212

$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 71156
$ws.Range("C51").Value = 691
$ws.Range("D51").Value = 46676
$ws.Range("E51").Value = 22549
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 1931

$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 71083
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 29253
$ws.Range("E52").Value = 40689
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 1141

$ws.Range("A58").Value = "Singapur"
$ws.Range("B58").Value = 57654
$ws.Range("C58").Value = 15
$ws.Range("D58").Value = 57333
$ws.Range("E58").Value = 294
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 27

$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("B60").Value = 53834
$ws.Range("C60").Value = 559
$ws.Range("D60").Value = 50230
$ws.Range("E60").Value = 3158
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 446

$ws.Range("A93").Value = "Noruega"
$ws.Range("B93").Value = 13347
$ws.Range("C93").Value = 70
$ws.Range("D93").Value = 10371
$ws.Range("E93").Value = 2706
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 270

$ws.Range("A94").Value = "Tunez"
$ws.Range("B94").Value = 13305
$ws.Range("C94").Value = 826
$ws.Range("D94").Value = 5032
$ws.Range("E94").Value = 8093
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 6
$ws.Range("H94").Value = 180

$ws.Range("A125").Value = "Eslovenia"
$ws.Range("B125").Value = 4816
$ws.Range("C125").Value = 122
$ws.Range("D125").Value = 3245
$ws.Range("E125").Value = 1426
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 145

$ws.Range("A126").Value = "Republica de Africa Central"
$ws.Range("B126").Value = 4802
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 1830
$ws.Range("E126").Value = 2910
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 62

$ws.Range("A127").Value = "Ruanda"
$ws.Range("B127").Value = 4779
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 2995
$ws.Range("E127").Value = 1757
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 27

$ws.Range("A128").Value = "Surinam"
$ws.Range("B128").Value = 4779
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 4560
$ws.Range("E128").Value = 118
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 101

$ws.Range("A142").Value = "Sri Lanka"
$ws.Range("B142").Value = 3327
$ws.Range("C142").Value = 3
$ws.Range("D142").Value = 3142
$ws.Range("E142").Value = 172
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 13

$ws.Range("A179").Value = "Islas Feroe"
$ws.Range("B179").Value = 455
$ws.Range("C179").Value = 4
$ws.Range("D179").Value = 416
$ws.Range("E179").Value = 39
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

$ws.Range("A181").Value = "Mauricio"
$ws.Range("B181").Value = 367
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 343
$ws.Range("E181").Value = 14
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 10

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
